$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 00:35"

# Row 4 - Estados Unidos: updated totals
$ws.Range("B4").Value = 1684699
$ws.Range("C4").Value = 17871
$ws.Range("D4").Value = 451567
$ws.Range("E4").Value = 1133874
$ws.Range("G4").Value = 575
$ws.Range("H4").Value = 99258

# Row 55 - Noruega: updated totals (D unchanged)
$ws.Range("B55").Value = 8352
$ws.Range("C55").Value = 6
$ws.Range("E55").Value = 390

# Row 58 now becomes Nigeria (new data), row 59 becomes Oman (old row-58 Oman data,
# shifted down because Nigeria's updated count now exceeds Oman's)
$ws.Range("A58").Value = "Nigeria"
$ws.Range("B58").Value = 7839
$ws.Range("C58").Value = 313
$ws.Range("D58").Value = 2263
$ws.Range("E58").Value = 5350
$ws.Range("G58").Value = 5
$ws.Range("H58").Value = 226

$ws.Range("A59").Value = "Oman"
$ws.Range("B59").Value = 7770
$ws.Range("C59").Value = 513
$ws.Range("D59").Value = 1933
$ws.Range("E59").Value = 5800
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 37

# Row 160 - Guadalupe: updated totals (D unchanged)
$ws.Range("B160").Value = 161
$ws.Range("C160").Value = 5
$ws.Range("E160").Value = 32
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 14
